# Apply the "Mass Balance Test" content edits:
#  - Tidy up several TEST STEPS cells in column B (remove trailing periods / "-> Save" suffixes)
#  - Move the active window / selection to H12 (previously G26)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B text corrections (order matches the original author's edit sequence)
$ws.Range("B12").Value = "Add new station"
$ws.Range("B13").Value = "Edit a station"
$ws.Range("B14").Value = "Delete station"
$ws.Range("B9").Value  = 'Program contains all MAB listed variables, buttons, and text on a tab labeled "Mass Balance"'
$ws.Range("B8").Value  = "Program recieves xml input for the pointmass variables related to the aircraft"
$ws.Range("B7").Value  = "Program recieves xml input for the units of moment of inertia variables"
$ws.Range("B6").Value  = "Program recieves xml input for six variables relating to the moment of inertia of the aircraft"
$ws.Range("B5").Value  = "Program recieves xml input for the units of location variable"
$ws.Range("B4").Value  = "Program recieves xml input for the location of the center of mass of the aircraft"
$ws.Range("B3").Value  = "Program recieves xml input for the units of mass variable"
$ws.Range("B2").Value  = "Program recieves xml input for the mass of the aircraft"

# Update the active cell selection shown when the workbook is reopened
$ws.Range("H12").Select()

$wb.Save()
